# Generate Report for Handoff
# Adds two newly-handed-off files (4a773841-... and 681f826b-...) to the
# localization status report. Each of the three sheets (Overview, zh-cn,
# de-de) gains two new rows, and the ".localization-config" row (which was
# previously the last row on every sheet) is pushed down below them.

$wb = $excel.ActiveWorkbook

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/ce9cd13ec15b414a8f94b9baac3f36df277b1619"
$zhHoBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd39251974e34bd23e9013bdedc177d33c903953/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHoBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c7f5b416ce5bca84fd3a219269126d95c98fa84e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$file1 = "4a773841-d93f-454e-8e9e-9230eb64a205"
$file2 = "681f826b-0c2a-464f-91e2-90e55434c729"
$xlf1Hash = "d370ec9ce554d607679e77b57a8e90107b0df2a4"
$xlf2Hash = "fb8e2afd8b044976739b16be46add8185604853b"

$zhDate = "2016-03-03 06:34:42"
$deDate = "2016-03-03 06:34:54"
$noDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"  (columns: File Name | zh-cn | de-de)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Range("B6").Value = "Not to be localized"
$ws1.Range("C6").Value = "Not to be localized"

# Rebuild hyperlinks in column A (existing ones are kept identical, the
# ones that shift down / change target are re-created with the right text)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "$repoBase/e2e/1a461fb9-6d41-4917-89a6-16a06b6dbc3b.md", "", "", "1a461fb9-6d41-4917-89a6-16a06b6dbc3b.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "$repoBase/e2e/b75d4aa1-3249-447d-9b8e-786f07c45f7c.md", "", "", "b75d4aa1-3249-447d-9b8e-786f07c45f7c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "$repoBase/e2e/$file1.md", "", "", "$file1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "$repoBase/e2e/$file2.md", "", "", "$file2.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "$repoBase/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("D4").Value = $zhDate
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $noDate
$ws2.Range("H4").Value = "Include"

$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("D5").Value = $zhDate
$ws2.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G5").Value = $noDate
$ws2.Range("H5").Value = "Include"

$ws2.Range("B6").Value = "Not to be localized"
$ws2.Range("D6").Value = $noDate
$ws2.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G6").Value = $noDate
$ws2.Range("H6").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$repoBase/e2e/1a461fb9-6d41-4917-89a6-16a06b6dbc3b.md", "", "", "1a461fb9-6d41-4917-89a6-16a06b6dbc3b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "$zhHoBase/1a461fb9-6d41-4917-89a6-16a06b6dbc3b.31e13ef5270c5ef7169d0983b1f70bc40516706e.zh-cn.xlf", "", "", "1a461fb9-6d41-4917-89a6-16a06b6dbc3b.31e13ef5270c5ef7169d0983b1f70bc40516706e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "$repoBase/e2e/b75d4aa1-3249-447d-9b8e-786f07c45f7c.md", "", "", "b75d4aa1-3249-447d-9b8e-786f07c45f7c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "$zhHoBase/b75d4aa1-3249-447d-9b8e-786f07c45f7c.6dac9fea1103e5f74b9f83658df723f52268386a.zh-cn.xlf", "", "", "b75d4aa1-3249-447d-9b8e-786f07c45f7c.6dac9fea1103e5f74b9f83658df723f52268386a.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "$repoBase/e2e/$file1.md", "", "", "$file1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "$zhHoBase/$file1.$xlf1Hash.zh-cn.xlf", "", "", "$file1.$xlf1Hash.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "$repoBase/e2e/$file2.md", "", "", "$file2.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "$zhHoBase/$file2.$xlf2Hash.zh-cn.xlf", "", "", "$file2.$xlf2Hash.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "$repoBase/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("D4").Value = $deDate
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $noDate
$ws3.Range("H4").Value = "Include"

$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("D5").Value = $deDate
$ws3.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G5").Value = $noDate
$ws3.Range("H5").Value = "Include"

$ws3.Range("B6").Value = "Not to be localized"
$ws3.Range("D6").Value = $noDate
$ws3.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G6").Value = $noDate
$ws3.Range("H6").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$repoBase/e2e/1a461fb9-6d41-4917-89a6-16a06b6dbc3b.md", "", "", "1a461fb9-6d41-4917-89a6-16a06b6dbc3b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "$deHoBase/1a461fb9-6d41-4917-89a6-16a06b6dbc3b.31e13ef5270c5ef7169d0983b1f70bc40516706e.de-de.xlf", "", "", "1a461fb9-6d41-4917-89a6-16a06b6dbc3b.31e13ef5270c5ef7169d0983b1f70bc40516706e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "$repoBase/e2e/b75d4aa1-3249-447d-9b8e-786f07c45f7c.md", "", "", "b75d4aa1-3249-447d-9b8e-786f07c45f7c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "$deHoBase/b75d4aa1-3249-447d-9b8e-786f07c45f7c.6dac9fea1103e5f74b9f83658df723f52268386a.de-de.xlf", "", "", "b75d4aa1-3249-447d-9b8e-786f07c45f7c.6dac9fea1103e5f74b9f83658df723f52268386a.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "$repoBase/e2e/$file1.md", "", "", "$file1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "$deHoBase/$file1.$xlf1Hash.de-de.xlf", "", "", "$file1.$xlf1Hash.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "$repoBase/e2e/$file2.md", "", "", "$file2.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "$deHoBase/$file2.$xlf2Hash.de-de.xlf", "", "", "$file2.$xlf2Hash.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "$repoBase/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Report regenerated for handoff: added $file1 and $file2 to Overview, zh-cn and de-de sheets."
